$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price/volume cells so values keep their literal
# string representation (trailing zeros, % signs) instead of being
# auto-converted to numbers/percentages by Excel.
$textCells = @("D2","E2","D3","E3","D4","E4","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","E17","D18","E18","D19","E19","D20","E20","E21","D22","E22","D23","E23","E24","D26","E26","E27","E28","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","E46","D47","E47","D48","E48","E49","E50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "246.48"
$ws.Range("E2").Value = "0.87%"
$ws.Range("D3").Value = "29.45"
$ws.Range("E3").Value = "7.37%"
$ws.Range("D4").Value = "5.198"
$ws.Range("E4").Value = "3.02%"
$ws.Range("E5").Value = "0.65%"
$ws.Range("D6").Value = "6.574"
$ws.Range("E6").Value = "1.48%"
$ws.Range("D7").Value = "0.8587"
$ws.Range("E7").Value = "4.50%"
$ws.Range("D8").Value = "0.8780"
$ws.Range("E8").Value = "4.69%"
$ws.Range("D9").Value = "0.1364"
$ws.Range("E9").Value = "2.89%"
$ws.Range("D10").Value = "0.07073"
$ws.Range("E10").Value = "2.02%"
$ws.Range("D11").Value = "0.02877"
$ws.Range("E11").Value = "0.47%"
$ws.Range("D12").Value = "0.09388"
$ws.Range("E12").Value = "-0.03%"
$ws.Range("D13").Value = "0.001512"
$ws.Range("E13").Value = "-0.21%"
$ws.Range("D14").Value = "0.04157"
$ws.Range("E14").Value = "0.56%"
$ws.Range("D15").Value = "0.0005997"
$ws.Range("E15").Value = "-94.00%"
$ws.Range("D16").Value = "0.006149"
$ws.Range("E16").Value = "0.07%"
$ws.Range("E17").Value = "5,108.04%"
$ws.Range("D18").Value = "3.482"
$ws.Range("E18").Value = "-0.76%"
$ws.Range("D19").Value = "3.063"
$ws.Range("E19").Value = "2.07%"
$ws.Range("D20").Value = "2.261"
$ws.Range("E20").Value = "-1.99%"
$ws.Range("E21").Value = "2.06%"
$ws.Range("D22").Value = "0.03296"
$ws.Range("E22").Value = "4.31%"
$ws.Range("D23").Value = "0.1300"
$ws.Range("E23").Value = "0.67%"
$ws.Range("E24").Value = "-2.59%"
$ws.Range("D26").Value = "0.005055"
$ws.Range("E26").Value = "30.50%"
$ws.Range("E27").Value = "0.03%"
$ws.Range("E28").Value = "23.43%"
$ws.Range("D40").Value = "0.03747"
$ws.Range("E40").Value = "1.36%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.005678"
$ws.Range("E41").Value = "-6.61%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1072"
$ws.Range("E42").Value = "1.86%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002099"
$ws.Range("E43").Value = "-9.03%"
$ws.Range("D44").Value = "0.009857"
$ws.Range("E44").Value = "5.70%"
$ws.Range("D45").Value = "0.00005130"
$ws.Range("E45").Value = "-1.24%"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("D47").Value = "0.07096"
$ws.Range("E47").Value = "-30.08%"
$ws.Range("D48").Value = "0.002591"
$ws.Range("E48").Value = "1.17%"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("E50").Value = "-0.04%"
